$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" "27.482.50"
Set-TextValue "E2" "  +1.61%  "
Set-TextValue "D3" "1.567.09"
Set-TextValue "D4" "0.992"
Set-TextValue "E4" "  -1.36%  "
Set-TextValue "D5" "211.20"
Set-TextValue "E5" "  +1.17%  "
Set-TextValue "E6" "  -0.20%  "
Set-TextValue "E7" "  -1.34%  "
Set-TextValue "D8" "22.64"
Set-TextValue "E8" "  +2.33%  "
Set-TextValue "E9" "  +0.38%  "
Set-TextValue "D10" "0.0595"
Set-TextValue "D11" "0.0871"
Set-TextValue "E11" "  +1.42%  "
Set-TextValue "D12" "1.789.06"
Set-TextValue "E12" "  +0.02%  "
Set-TextValue "D13" "1.563.92"
Set-TextValue "E13" "  -0.13%  "
Set-TextValue "D14" "3.75"
Set-TextValue "E14" "  -1.14%  "
Set-TextValue "D15" "0.519"
Set-TextValue "E15" "  -0.40%  "
Set-TextValue "D16" "27.452.02"
Set-TextValue "E16" "  +1.51%  "
Set-TextValue "D17" "62.36"
Set-TextValue "E17" "  +0.72%  "
Set-TextValue "D18" "225.29"
Set-TextValue "E18" "  +4.09%  "
Set-TextValue "D19" "7.49"
Set-TextValue "E19" "  +0.82%  "
Set-TextValue "E20" "  -0.57%  "
Set-TextValue "D21" "0.993"
Set-TextValue "E21" "  -1.27%  "
Set-TextValue "D22" "4.11"
Set-TextValue "D23" "9.38"
Set-TextValue "E23" "  +1.48%  "
Set-TextValue "E24" "  +0.58%  "
Set-TextValue "D25" "150.04"
Set-TextValue "E25" "  -2.60%  "
Set-TextValue "B26" "EthereumClassic"
Set-TextValue "C26" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D26" "15.14"
Set-TextValue "E26" "  +0.39%  "
Set-TextValue "D27" "6.60"
Set-TextValue "B28" "Stellar"
Set-TextValue "C28" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D28" "0.107"
Set-TextValue "E28" "  +1.59%  "
Set-TextValue "D29" "0.992"
Set-TextValue "E29" "  -1.34%  "
Set-TextValue "E30" "  +0.77%  "
Set-TextValue "D31" "0.0471"
Set-TextValue "E31" "  -0.97%  "
Set-TextValue "E32" "  -0.27%  "
Set-TextValue "D33" "1.448.59"
Set-TextValue "E33" "  +1.54%  "
Set-TextValue "E34" "  -2.03%  "
Set-TextValue "E35" "  +3.25%  "
Set-TextValue "E36" "  -0.95%  "
Set-TextValue "E37" "  -0.77%  "
Set-TextValue "E39" "  +1.11%  "
Set-TextValue "D40" "0.813"
Set-TextValue "E40" "  -0.06%  "
Set-TextValue "B41" "MXToken"
Set-TextValue "C41" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D41" "2.37"
Set-TextValue "E41" "  +1.54%  "
Set-TextValue "B42" "FraxShare"
Set-TextValue "C42" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D42" "5.74"
Set-TextValue "E42" "  -1.24%  "
Set-TextValue "E43" "  -1.34%  "
Set-TextValue "D44" "1.86"
Set-TextValue "E44" "  +6.31%  "
Set-TextValue "E45" "  -2.93%  "
Set-TextValue "D46" "64.23"
Set-TextValue "E46" "  -1.00%  "
Set-TextValue "D47" "1.701.78"
Set-TextValue "E47" "  -0.07%  "
Set-TextValue "D48" "86.83"
Set-TextValue "E48" "  +0.14%  "
Set-TextValue "E49" "  +0.54%  "
Set-TextValue "D50" "0.0526"
Set-TextValue "E50" "  +1.18%  "
Set-TextValue "E51" "  -2.00%  "
